$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ligand symbol (B) and Receptor symbol (C) columns for data rows 2-9
$ws.Range("B2:B9").Value = "Efna3"
$ws.Range("C2:C9").Value = "Epha2"

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3193606666666667
$ws.Range("H2").Value = 0.958082
$ws.Range("I2").Value = 0.7979421849584948
$ws.Range("J2").Value = 0.7979421849584948
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.163974
$ws.Range("N2").Value = 57.491922
$ws.Range("O2").Value = 0.6845732287637933
$ws.Range("P2").Value = 0.6845732287637933
$ws.Range("Q2").Value = 6.120219512622667
$ws.Range("R2").Value = 55.081975613604
$ws.Range("S2").Value = 0.5462498579238727
$ws.Range("T2").Value = 0.5462498579238727

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3193606666666667
$ws.Range("H3").Value = 0.958082
$ws.Range("I3").Value = 0.7979421849584948
$ws.Range("J3").Value = 0.7979421849584948
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7651789999999999
$ws.Range("N3").Value = 2.295537
$ws.Range("O3").Value = 0.02733363438148322
$ws.Range("P3").Value = 0.02733363438148323
$ws.Range("Q3").Value = 0.2443680755593333
$ws.Range("R3").Value = 2.199312680034
$ws.Range("S3").Value = 0.02181065994121736
$ws.Range("T3").Value = 0.02181065994121736

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3193606666666667
$ws.Range("H4").Value = 0.958082
$ws.Range("I4").Value = 0.7979421849584948
$ws.Range("J4").Value = 0.7979421849584948
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.880893333333333
$ws.Range("N4").Value = 23.64268
$ws.Range("O4").Value = 0.281520346184098
$ws.Range("P4").Value = 0.281520346184098
$ws.Range("Q4").Value = 2.516847348862222
$ws.Range("R4").Value = 22.65162613976
$ws.Range("S4").Value = 0.224636960144411
$ws.Range("T4").Value = 0.224636960144411

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3193606666666667
$ws.Range("H5").Value = 0.958082
$ws.Range("I5").Value = 0.7979421849584948
$ws.Range("J5").Value = 0.7979421849584948
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.183999
$ws.Range("N5").Value = 0.551997
$ws.Range("O5").Value = 0.006572790670625477
$ws.Range("P5").Value = 0.006572790670625476
$ws.Range("Q5").Value = 0.058762043306
$ws.Range("R5").Value = 0.5288583897539999
$ws.Range("S5").Value = 0.005244706948993704
$ws.Range("T5").Value = 0.005244706948993703

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08086966666666666
$ws.Range("H6").Value = 0.242609
$ws.Range("I6").Value = 0.2020578150415052
$ws.Range("J6").Value = 0.2020578150415053
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.163974
$ws.Range("N6").Value = 57.491922
$ws.Range("O6").Value = 0.6845732287637933
$ws.Range("P6").Value = 0.6845732287637933
$ws.Range("Q6").Value = 1.549784189388667
$ws.Range("R6").Value = 13.948057704498
$ws.Range("S6").Value = 0.1383233708399206
$ws.Range("T6").Value = 0.1383233708399206

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08086966666666666
$ws.Range("H7").Value = 0.242609
$ws.Range("I7").Value = 0.2020578150415052
$ws.Range("J7").Value = 0.2020578150415053
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7651789999999999
$ws.Range("N7").Value = 2.295537
$ws.Range("O7").Value = 0.02733363438148322
$ws.Range("P7").Value = 0.02733363438148323
$ws.Range("Q7").Value = 0.06187977067033332
$ws.Range("R7").Value = 0.556917936033
$ws.Range("S7").Value = 0.005522974440265865
$ws.Range("T7").Value = 0.005522974440265867

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.08086966666666666
$ws.Range("H8").Value = 0.242609
$ws.Range("I8").Value = 0.2020578150415052
$ws.Range("J8").Value = 0.2020578150415053
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.880893333333333
$ws.Range("N8").Value = 23.64268
$ws.Range("O8").Value = 0.281520346184098
$ws.Range("P8").Value = 0.281520346184098
$ws.Range("Q8").Value = 0.6373252169022221
$ws.Range("R8").Value = 5.73592695212
$ws.Range("S8").Value = 0.056883386039687
$ws.Range("T8").Value = 0.056883386039687

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.08086966666666666
$ws.Range("H9").Value = 0.242609
$ws.Range("I9").Value = 0.2020578150415052
$ws.Range("J9").Value = 0.2020578150415053
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.183999
$ws.Range("N9").Value = 0.551997
$ws.Range("O9").Value = 0.006572790670625477
$ws.Range("P9").Value = 0.006572790670625476
$ws.Range("Q9").Value = 0.014879937797
$ws.Range("R9").Value = 0.133919440173
$ws.Range("S9").Value = 0.001328083721631774
$ws.Range("T9").Value = 0.001328083721631774

# Remove now-obsolete rows 10-13 (data reduced from 13 to 9 rows)
$ws.Rows("10:13").Delete()
